$d = $word.ActiveDocument

$d.Content.Find.Execute("594×9=", $true, $false, $false, $false, $false, $true, 1, $false, "474×4=", 2) | Out-Null
$d.Content.Find.Execute("426×6=", $true, $false, $false, $false, $false, $true, 1, $false, "235×3=", 2) | Out-Null
$d.Content.Find.Execute("674×8=", $true, $false, $false, $false, $false, $true, 1, $false, "589×7=", 2) | Out-Null
$d.Content.Find.Execute("396×2=", $true, $false, $false, $false, $false, $true, 1, $false, "943×5=", 2) | Out-Null
$d.Content.Find.Execute("890×3=", $true, $false, $false, $false, $false, $true, 1, $false, "135×9=", 2) | Out-Null
$d.Content.Find.Execute("707×8=", $true, $false, $false, $false, $false, $true, 1, $false, "154×3=", 2) | Out-Null
$d.Content.Find.Execute("202×5=", $true, $false, $false, $false, $false, $true, 1, $false, "479×4=", 2) | Out-Null
$d.Content.Find.Execute("396×5=", $true, $false, $false, $false, $false, $true, 1, $false, "910×6=", 2) | Out-Null
$d.Content.Find.Execute("781×8=", $true, $false, $false, $false, $false, $true, 1, $false, "876×9=", 2) | Out-Null
$d.Content.Find.Execute("680×4=", $true, $false, $false, $false, $false, $true, 1, $false, "628×5=", 2) | Out-Null
$d.Content.Find.Execute("799×8=", $true, $false, $false, $false, $false, $true, 1, $false, "485×6=", 2) | Out-Null
$d.Content.Find.Execute("953×8=", $true, $false, $false, $false, $false, $true, 1, $false, "893×4=", 2) | Out-Null
$d.Content.Find.Execute("669×4=", $true, $false, $false, $false, $false, $true, 1, $false, "299×6=", 2) | Out-Null
$d.Content.Find.Execute("413×8=", $true, $false, $false, $false, $false, $true, 1, $false, "536×6=", 2) | Out-Null
$d.Content.Find.Execute("845×7=", $true, $false, $false, $false, $false, $true, 1, $false, "250×4=", 2) | Out-Null
$d.Content.Find.Execute("430×2=", $true, $false, $false, $false, $false, $true, 1, $false, "527×3=", 2) | Out-Null
$d.Content.Find.Execute("358×2=", $true, $false, $false, $false, $false, $true, 1, $false, "340×3=", 2) | Out-Null
$d.Content.Find.Execute("525×7=", $true, $false, $false, $false, $false, $true, 1, $false, "499×5=", 2) | Out-Null
$d.Content.Find.Execute("173×4=", $true, $false, $false, $false, $false, $true, 1, $false, "178×8=", 2) | Out-Null
$d.Content.Find.Execute("802×2=", $true, $false, $false, $false, $false, $true, 1, $false, "300×2=", 2) | Out-Null
$d.Content.Find.Execute("579×5=", $true, $false, $false, $false, $false, $true, 1, $false, "157×2=", 2) | Out-Null
$d.Content.Find.Execute("261×8=", $true, $false, $false, $false, $false, $true, 1, $false, "392×5=", 2) | Out-Null
$d.Content.Find.Execute("409×4=", $true, $false, $false, $false, $false, $true, 1, $false, "971×3=", 2) | Out-Null
$d.Content.Find.Execute("495×3=", $true, $false, $false, $false, $false, $true, 1, $false, "155×7=", 2) | Out-Null
$d.Content.Find.Execute("708×3=", $true, $false, $false, $false, $false, $true, 1, $false, "336×8=", 2) | Out-Null
